$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7750.25
$ws.Range("I19").Value = 20002
$ws.Range("J19").Value = 3666.3333
$ws.Range("K19").Value = 20002
$ws.Range("L19").Value = 3666.3333
$ws.Range("M19").Value = -19827
$ws.Range("N19").Value = -4016.3333
$ws.Range("H58").Value = 1968.8
$ws.Range("I58").Value = 146.11111
$ws.Range("J58").Value = 4702.8335
$ws.Range("K58").Value = 438.33333
$ws.Range("L58").Value = 14108.5005
$ws.Range("M58").Value = -288.33333
$ws.Range("N58").Value = -14408.5005
$ws.Range("H98").Value = 2281.5386
$ws.Range("I98").Value = 2345.1714
$ws.Range("J98").Value = 1724.75
$ws.Range("K98").Value = 2345.1714
$ws.Range("L98").Value = 1724.75
$ws.Range("M98").Value = -847.1714000000002
$ws.Range("N98").Value = -4720.75
$ws.Range("H122").Value = 2281.5386
$ws.Range("I122").Value = 2345.1714
$ws.Range("J122").Value = 1724.75
$ws.Range("K122").Value = 7035.514200000001
$ws.Range("L122").Value = 5174.25
$ws.Range("M122").Value = -4585.514200000001
$ws.Range("N122").Value = -10074.25
$ws.Range("H127").Value = 701.4761999999999
$ws.Range("I127").Value = 407.23077
$ws.Range("J127").Value = 1179.625
$ws.Range("K127").Value = 1221.69231
$ws.Range("L127").Value = 3538.875
$ws.Range("M127").Value = 3738.30769
$ws.Range("N127").Value = -13458.875
$ws.Range("H132").Value = 2967.739
$ws.Range("I132").Value = 3083.3171
$ws.Range("K132").Value = 9249.951300000001
$ws.Range("M132").Value = -6719.951300000001
$ws.Range("H137").Value = 21740324
$ws.Range("I137").Value = 1023.2647
$ws.Range("J137").Value = 83335010
$ws.Range("K137").Value = 3069.7941
$ws.Range("L137").Value = 250005030
$ws.Range("M137").Value = -519.7941000000001
$ws.Range("N137").Value = -250010130
$ws.Range("H138").Value = 2323.6543
$ws.Range("I138").Value = 1802.6364
$ws.Range("J138").Value = 2943.2432
$ws.Range("K138").Value = 5407.9092
$ws.Range("L138").Value = 8829.729599999999
$ws.Range("M138").Value = -267.9092000000001
$ws.Range("N138").Value = -19109.7296
$ws.Range("H140").Value = 65000
$ws.Range("J140").Value = 65000
$ws.Range("L140").Value = 65000
$ws.Range("N140").Value = -75360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7855.9697
$ws.Range("I32").Value = 4662
$ws.Range("J32").Value = 17310.12
$ws.Range("K32").Value = 4662
$ws.Range("L32").Value = 17310.12
$ws.Range("M32").Value = -4375
$ws.Range("N32").Value = -17884.12
$ws.Range("H61").Value = 2416433.5
$ws.Range("I61").Value = 2416433.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2416433.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2416221.5
$ws.Range("N61").ClearContents()
$ws.Range("H88").Value = 2591.2666
$ws.Range("I88").Value = 2903.6667
$ws.Range("J88").Value = 2383
$ws.Range("K88").Value = 2903.6667
$ws.Range("L88").Value = 2383
$ws.Range("M88").Value = -2497.6667
$ws.Range("N88").Value = -3195
$ws.Range("H91").Value = 2591.2666
$ws.Range("I91").Value = 2903.6667
$ws.Range("J91").Value = 2383
$ws.Range("K91").Value = 2903.6667
$ws.Range("L91").Value = 2383
$ws.Range("M91").Value = -1499.6667
$ws.Range("N91").Value = -5191
$ws.Range("H132").Value = 8753426
$ws.Range("I132").Value = 10074120
$ws.Range("J132").Value = 168914
$ws.Range("K132").Value = 30222360
$ws.Range("L132").Value = 506742
$ws.Range("M132").Value = -30219830
$ws.Range("N132").Value = -511802
$ws.Range("H136").Value = 2416433.5
$ws.Range("I136").Value = 2416433.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7249300.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7246750.5
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2194.9148
$ws.Range("I86").Value = 1783.4814
$ws.Range("J86").Value = 2750.35
$ws.Range("K86").Value = 1783.4814
$ws.Range("L86").Value = 2750.35
$ws.Range("M86").Value = -660.4813999999999
$ws.Range("N86").Value = -4996.35
$ws.Range("H89").Value = 2194.9148
$ws.Range("I89").Value = 1783.4814
$ws.Range("J89").Value = 2750.35
$ws.Range("K89").Value = 8917.406999999999
$ws.Range("L89").Value = 13751.75
$ws.Range("M89").Value = -3301.406999999999
$ws.Range("N89").Value = -24983.75
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620
$ws.Range("H130").Value = 40780
$ws.Range("J130").Value = 40780
$ws.Range("L130").Value = 40780
$ws.Range("N130").Value = -50820
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7382398
$ws.Range("I31").Value = 1414.4546
$ws.Range("J31").Value = 61509610
$ws.Range("K31").Value = 1414.4546
$ws.Range("L31").Value = 61509610
$ws.Range("M31").Value = -1119.4546
$ws.Range("N31").Value = -61510200
$ws.Range("H34").Value = 7382398
$ws.Range("I34").Value = 1414.4546
$ws.Range("J34").Value = 61509610
$ws.Range("K34").Value = 1414.4546
$ws.Range("L34").Value = 61509610
$ws.Range("M34").Value = -1212.4546
$ws.Range("N34").Value = -61510014
$ws.Range("H62").Value = 2683.3333
$ws.Range("I62").Value = 2250
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 2250
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -1626
$ws.Range("N62").Value = -4148
$ws.Range("H65").Value = 2683.3333
$ws.Range("I65").Value = 2250
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 11250
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -8130
$ws.Range("N65").Value = -20740
$ws.Range("H99").Value = 3200
$ws.Range("I99").Value = 3200
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3200
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1702
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 3200
$ws.Range("I126").Value = 3200
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9600
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7130
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2021.238
$ws.Range("I132").Value = 1949.5518
$ws.Range("J132").Value = 2181.1538
$ws.Range("K132").Value = 5848.6554
$ws.Range("L132").Value = 6543.4614
$ws.Range("M132").Value = -3318.6554
$ws.Range("N132").Value = -11603.4614
$ws.Range("H134").Value = 1251.8718
$ws.Range("I134").Value = 1469.5186
$ws.Range("J134").Value = 762.1667
$ws.Range("K134").Value = 4408.5558
$ws.Range("L134").Value = 2286.5001
$ws.Range("M134").Value = -1873.5558
$ws.Range("N134").Value = -7356.5001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 300
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -770
$ws.Range("H105").Value = 8000000
$ws.Range("J105").Value = 8000000
$ws.Range("L105").Value = 24000000
$ws.Range("N105").Value = -24005242
$ws.Range("H131").Value = 3254.818
$ws.Range("J131").Value = 2320.647
$ws.Range("L131").Value = 6961.941
$ws.Range("N131").Value = -17041.941
$ws.Range("H132").Value = 200001800
$ws.Range("I132").Value = 500000260
$ws.Range("K132").Value = 4500002340
$ws.Range("M132").Value = -4499999810
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4975
$ws.Range("I70").Value = 6000
$ws.Range("J70").Value = 4633.3335
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 4633.3335
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -5173.3335
$ws.Range("H73").Value = 4975
$ws.Range("I73").Value = 6000
$ws.Range("J73").Value = 4633.3335
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 4633.3335
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -6505.3335
$ws.Range("H102").Value = 1754.1765
$ws.Range("I102").Value = 1663.2
$ws.Range("J102").Value = 1884.1428
$ws.Range("K102").Value = 1663.2
$ws.Range("L102").Value = 1884.1428
$ws.Range("M102").Value = -41.20000000000005
$ws.Range("N102").Value = -5128.1428
$ws.Range("H126").Value = 2052
$ws.Range("I126").Value = 1878
$ws.Range("K126").Value = 5634
$ws.Range("M126").Value = -3164
$ws.Range("H141").Value = 141885
$ws.Range("J141").Value = 141885
$ws.Range("L141").Value = 141885
$ws.Range("N141").Value = -152245
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2126.6667
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2380
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2380
$ws.Range("M40").Value = -1864
$ws.Range("N40").Value = -2652
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 68857.5
$ws.Range("J141").Value = 68857.5
$ws.Range("L141").Value = 68857.5
$ws.Range("N141").Value = -79217.5
